$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column in H1, matching the style of the existing
# header row (same bold/border/centered format used by B1:G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H8 with the "Save" values (plain numbers, no special style,
# matching the rest of the numeric data columns).
$saveValues = @(1, 0, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}

$excel.CutCopyMode = $false
